$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 5.3333335
$ws.Range("I11").Value = 5.3333335
$ws.Range("K11").Value = 5.3333335
$ws.Range("M11").Value = 134.6666665
$ws.Range("H51").Value = 5140.857
$ws.Range("J51").Value = 5165.1665
$ws.Range("L51").Value = 5165.1665
$ws.Range("N51").Value = -6133.1665
$ws.Range("H86").Value = 9960.267
$ws.Range("I86").Value = 4468.625
$ws.Range("K86").Value = 4468.625
$ws.Range("M86").Value = -3345.625
$ws.Range("H89").Value = 9960.267
$ws.Range("I89").Value = 4468.625
$ws.Range("K89").Value = 22343.125
$ws.Range("M89").Value = -16727.125
$ws.Range("H113").Value = 6246.25
$ws.Range("J113").Value = 8560
$ws.Range("L113").Value = 8560
$ws.Range("N113").Value = -15068
$ws.Range("H135").Value = 1364.375
$ws.Range("I135").Value = 1314.6666
$ws.Range("J135").Value = 1394.2
$ws.Range("K135").Value = 11831.9994
$ws.Range("L135").Value = 12547.8
$ws.Range("M135").Value = -9296.999400000001
$ws.Range("N135").Value = -17617.8
$ws.Range("H137").Value = 2137.6428
$ws.Range("I137").Value = 1629.8182
$ws.Range("K137").Value = 4889.4546
$ws.Range("M137").Value = -2339.4546
$ws.Range("H141").Value = 6709.3335
$ws.Range("I141").Value = 6709.3335
$ws.Range("K141").Value = 20128.0005
$ws.Range("M141").Value = -14948.0005

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 2000000
$ws.Range("I11").Value = 2000000
$ws.Range("K11").Value = 2000000
$ws.Range("M11").Value = -1999856
$ws.Range("H55").Value = 18333.334
$ws.Range("I55").Value = 5000
$ws.Range("J55").Value = 25000
$ws.Range("K55").Value = 5000
$ws.Range("L55").Value = 25000
$ws.Range("M55").Value = -4685
$ws.Range("N55").Value = -25630
$ws.Range("H61").Value = 3072
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 2038.1333
$ws.Range("I74").Value = 1469.6428
$ws.Range("K74").Value = 1469.6428
$ws.Range("M74").Value = -595.6428000000001
$ws.Range("H77").Value = 2038.1333
$ws.Range("I77").Value = 1469.6428
$ws.Range("K77").Value = 7348.214
$ws.Range("M77").Value = -2980.214
$ws.Range("H97").Value = 1135
$ws.Range("I97").Value = 1135
$ws.Range("K97").Value = 1135
$ws.Range("M97").Value = -639
$ws.Range("H102").Value = 2000.875
$ws.Range("I102").Value = 2000.875
$ws.Range("K102").Value = 2000.875
$ws.Range("M102").Value = -378.875
$ws.Range("H114").Value = 28900
$ws.Range("J114").Value = 28900
$ws.Range("L114").Value = 28900
$ws.Range("N114").Value = -37578
$ws.Range("H132").Value = 1362.75
$ws.Range("I132").Value = 1362.75
$ws.Range("K132").Value = 4088.25
$ws.Range("M132").Value = -1558.25
$ws.Range("H136").Value = 3072
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("H139").Value = 44999.332
$ws.Range("J139").Value = 44999.332
$ws.Range("L139").Value = 44999.332
$ws.Range("N139").Value = -55279.332
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1198.9524
$ws.Range("I94").Value = 956.2105
$ws.Range("K94").Value = 956.2105
$ws.Range("M94").Value = -505.2105
$ws.Range("H134").Value = 12148.529
$ws.Range("I134").Value = 8940.6875
$ws.Range("K134").Value = 26822.0625
$ws.Range("M134").Value = -24287.0625

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 16752.75
$ws.Range("I58").Value = 12665.667
$ws.Range("K58").Value = 12665.667
$ws.Range("M58").Value = -12462.667
$ws.Range("H92").Value = 56720
$ws.Range("J92").Value = 56720
$ws.Range("L92").Value = 56720
$ws.Range("N92").Value = -61712
$ws.Range("H136").Value = 16752.75
$ws.Range("I136").Value = 12665.667
$ws.Range("K136").Value = 37997.001
$ws.Range("M136").Value = -35447.001

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3722501.5
$ws.Range("I4").Value = 143452
$ws.Range("J4").Value = 5810280.5
$ws.Range("K4").Value = 430356
$ws.Range("L4").Value = 17430841.5
$ws.Range("M4").Value = -430244
$ws.Range("N4").Value = -17431065.5
$ws.Range("H8").Value = 1143.1428
$ws.Range("I8").Value = 1143.1428
$ws.Range("K8").Value = 3429.4284
$ws.Range("M8").Value = -3290.4284
$ws.Range("H38").Value = 30.384615
$ws.Range("I38").Value = 28
$ws.Range("J38").Value = 43.5
$ws.Range("K38").Value = 84
$ws.Range("L38").Value = 130.5
$ws.Range("M38").Value = 263
$ws.Range("N38").Value = -824.5
$ws.Range("H50").Value = 407.55554
$ws.Range("I50").Value = 483.57144
$ws.Range("K50").Value = 1450.71432
$ws.Range("M50").Value = -969.71432
$ws.Range("H53").Value = 407.55554
$ws.Range("I53").Value = 483.57144
$ws.Range("K53").Value = 1450.71432
$ws.Range("M53").Value = -969.71432
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H97").Value = 927.8461
$ws.Range("I97").Value = 1680
$ws.Range("J97").Value = 593.55554
$ws.Range("K97").Value = 5040
$ws.Range("L97").Value = 1780.66662
$ws.Range("M97").Value = -4544
$ws.Range("N97").Value = -2772.66662
$ws.Range("H104").Value = 17375
$ws.Range("J104").Value = 29500
$ws.Range("L104").Value = 88500
$ws.Range("N104").Value = -93742
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("H137").Value = 5230.778
$ws.Range("J137").Value = 8310.799999999999
$ws.Range("L137").Value = 24932.4
$ws.Range("N137").Value = -35132.39999999999

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 1566
$ws.Range("I22").Value = 1566
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1566
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1037
$ws.Range("N22").ClearContents()
$ws.Range("H70").Value = 5198.5713
$ws.Range("I70").Value = 5198.5713
$ws.Range("K70").Value = 5198.5713
$ws.Range("M70").Value = -4928.5713
$ws.Range("H73").Value = 5198.5713
$ws.Range("I73").Value = 5198.5713
$ws.Range("K73").Value = 5198.5713
$ws.Range("M73").Value = -4262.5713
$ws.Range("H92").Value = 4125
$ws.Range("J92").Value = 4125
$ws.Range("L92").Value = 4125
$ws.Range("N92").Value = -7869
$ws.Range("H114").Value = 98796.60000000001
$ws.Range("J114").Value = 98796.60000000001
$ws.Range("L114").Value = 98796.60000000001
$ws.Range("N114").Value = -107474.6
$ws.Range("H132").Value = 933
$ws.Range("I132").Value = 933
$ws.Range("K132").Value = 2799
$ws.Range("M132").Value = -269
$ws.Range("H134").Value = 105332.664
$ws.Range("J134").Value = 105332.664
$ws.Range("L134").Value = 315997.992
$ws.Range("N134").Value = -321067.992

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4999.5
$ws.Range("J68").Value = 5999
$ws.Range("L68").Value = 5999
$ws.Range("N68").Value = -7497
$ws.Range("H71").Value = 4999.5
$ws.Range("J71").Value = 5999
$ws.Range("L71").Value = 29995
$ws.Range("N71").Value = -37483
$ws.Range("H122").Value = 4010.25
$ws.Range("I122").Value = 3632.5
$ws.Range("K122").Value = 10897.5
$ws.Range("M122").Value = -8447.5
$ws.Range("H136").Value = 3668
$ws.Range("J136").Value = 4500
$ws.Range("L136").Value = 13500
$ws.Range("N136").Value = -18600

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 25000
$ws.Range("J28").Value = 25000
$ws.Range("L28").Value = 25000
$ws.Range("N28").Value = -25696
$ws.Range("H94").Value = 25999.5
$ws.Range("J94").Value = 25999.5
$ws.Range("L94").Value = 25999.5
$ws.Range("N94").Value = -27801.5
$ws.Range("H103").Value = 50000
$ws.Range("J103").Value = 50000
$ws.Range("L103").Value = 50000
$ws.Range("N103").Value = -52344
$ws.Range("H105").Value = 20615
$ws.Range("J105").Value = 20615
$ws.Range("L105").Value = 20615
$ws.Range("N105").Value = -27603
$ws.Range("H140").Value = 86000
$ws.Range("J140").Value = 86000
$ws.Range("L140").Value = 86000
$ws.Range("N140").Value = -96360
